$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (date moved from 03-11 to 03-12)
$ws.Name = "Through 2022-03-12"

# Update the "March (through 03-11)" label in column A
$ws.Range("A4").Value = "March (through 03-12)"

# Update March row (row 4) values
$ws.Range("C4").Value = 16
$ws.Range("D4").Value = 24
$ws.Range("E4").Value = 23
$ws.Range("G4").Value = 23
$ws.Range("H4").Value = 34
$ws.Range("I4").Value = 56

# Update Total row (row 5) values
$ws.Range("C5").Value = 103
$ws.Range("D5").Value = 155
$ws.Range("E5").Value = 160
$ws.Range("G5").Value = 164
$ws.Range("H5").Value = 376
$ws.Range("I5").Value = 356
